$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": row 3 (b.md) moves from "Handed back: in sync with en-US"
# to "Ready for handoff", with a refreshed generation timestamp.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-17 06:33:03"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": row 3 (b.md) gets a new handback file, refreshed handback
# datetime, status flips to "Ready for handoff", Content Duplicate flips to
# False, and an error message is recorded for the stale handback.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-17 06:32:56"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6e35fa964110a17e7d2166f5c37226707d7e4369/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5920f454b55a54dcfd2201b736a32ce0b4f8d50d/e2e/b.md."
$wsZhCn.Columns("P").ColumnWidth = 39.14

# ---------------------------------------------------------------------------
# Sheet "de-de": same kind of refresh as zh-cn, but with the de-de handback
# file name and its own refreshed handback datetime.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-17 06:33:03"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6e35fa964110a17e7d2166f5c37226707d7e4369/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5920f454b55a54dcfd2201b736a32ce0b4f8d50d/e2e/b.md."
$wsDeDe.Columns("P").ColumnWidth = 39.14
